# add cpu % metric to data-set
# Appends a second data row (Date / Hour / CPU Utilization) below the
# existing header row on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Date -----------------------------------------------------
# Force text storage so "2021-03-30" isn't auto-parsed into a date serial,
# then drop back to the workbook's default "Normal" style so no stray
# number-format override is left behind on the cell.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2021-03-30"
$ws.Range("A2").Style = "Normal"

# --- Row 2: Hour -------------------------------------------------------
# Same trick so "21:28:27" stays literal text instead of becoming a time.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "21:28:27"
$ws.Range("B2").Style = "Normal"

# --- Row 2: CPU Utilization --------------------------------------------
$ws.Range("C2").Value = 0.508474576271109
